$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 42-72 (column C only decremented by 1) ---
$colCOnly = @{
    42 = 35; 43 = 39; 44 = 48; 45 = 52; 46 = 85; 47 = 104; 48 = 125; 49 = 159;
    50 = 189; 51 = 227; 52 = 299; 53 = 392; 54 = 510; 55 = 646; 56 = 815;
    57 = 1021; 58 = 1270; 59 = 1618; 60 = 2040; 61 = 2463; 62 = 2991;
    63 = 3425; 64 = 3905; 65 = 4437; 66 = 4983; 67 = 5712; 68 = 6407;
    69 = 7126; 70 = 7725; 71 = 8153; 72 = 8732
}
foreach ($r in $colCOnly.Keys) {
    $ws.Cells.Item($r, 3).Value = $colCOnly[$r]
}

# --- Update rows 73-82 (columns B and C decremented by 1) ---
$colBC = @{
    73 = @(130586, 9184);
    74 = @(137513, 9563);
    75 = @(143382, 9901);
    76 = @(149130, 10246);
    77 = @(156457, 10607);
    78 = @(163149, 10952);
    79 = @(173769, 11510);
    80 = @(184728, 11954);
    81 = @(197391, 12361);
    82 = @(206841, 12673)
}
foreach ($r in $colBC.Keys) {
    $vals = $colBC[$r]
    $ws.Cells.Item($r, 2).Value = $vals[0]
    $ws.Cells.Item($r, 3).Value = $vals[1]
}

# --- Update rows 83-88 (columns B, C, D, E, F all change) ---
$rows83to88 = @{
    83 = @(219624, 12975, 698, 180, 135);
    84 = @(231189, 13277, 657, 163, 124);
    85 = @(241820, 13588, 632, 168, 114);
    86 = @(253106, 13877, 610, 153, 114);
    87 = @(268245, 14180, 608, 150, 119);
    88 = @(281733, 14473, 538, 141, 116)
}
foreach ($r in $rows83to88.Keys) {
    $vals = $rows83to88[$r]
    $ws.Cells.Item($r, 2).Value = $vals[0]
    $ws.Cells.Item($r, 3).Value = $vals[1]
    $ws.Cells.Item($r, 4).Value = $vals[2]
    $ws.Cells.Item($r, 5).Value = $vals[3]
    $ws.Cells.Item($r, 6).Value = $vals[4]
}

# --- Update row 89 (columns C, D, E, F change; B unchanged) ---
$ws.Cells.Item(89, 3).Value = 14703
$ws.Cells.Item(89, 4).Value = 523
$ws.Cells.Item(89, 5).Value = 143
$ws.Cells.Item(89, 6).Value = 106

# --- Add new row 90 ---
$ws.Cells.Item(90, 1).Value = 43944
$ws.Cells.Item(90, 2).Value = 309188
$ws.Cells.Item(90, 3).Value = 14985
$ws.Cells.Item(90, 4).Value = 484
$ws.Cells.Item(90, 5).Value = 140
$ws.Cells.Item(90, 6).Value = 110
$ws.Cells.Item(90, 7).Value = 193

# Copy style from row 89 to row 90 so formatting matches the rest of the table
$ws.Range("A89:G89").Copy() | Out-Null
$ws.Range("A90:G90").PasteSpecial(-4122) | Out-Null

# --- Update selection ---
$ws.Range("J9").Select() | Out-Null
